$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New material labels (column A) and roughness values (column B) for rows 2-18
$materials = @(
    "cuivre, plomb, laiton, inox",
    "PVC",
    "acier inox",
    "tube acier du commerce",
    "acier étiré",
    "acier soudé",
    "acier galvanisé",
    "acier rouillé",
    "fonte neuve",
    "fonte usagée",
    "fonte incrustée",
    "tôle ou fonte asphaltée",
    "ciment bien lissé",
    "béton ordinaire",
    "béton grossier",
    "bois bien raboté",
    "bois ordinaire"
)

$values = @(
    0.001,
    0.0015,
    0.015,
    0.07,
    0.015,
    0.045,
    0.15,
    0.55,
    0.5,
    1.2,
    2,
    0.012,
    0.3,
    1,
    5,
    5,
    1
)

# Set "fonte usagée" (row 11, index 9) before "fonte neuve" (row 10, index 8)
# so that the shared-strings table gets the same insertion order as the
# target file (fonte usagée ends up with a lower shared-string index than
# fonte neuve, even though it is displayed one row below it).
$order = @(0, 1, 2, 3, 4, 5, 6, 7, 9, 8, 10, 11, 12, 13, 14, 15, 16)

foreach ($i in $order) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $materials[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Apply scientific number format to the PVC roughness value (B3)
$ws.Range("B3").NumberFormat = "0.00E+00"

# Update the selected cell like after the edits were made
$ws.Range("B19").Select()

$wb.Save()
